# closing dates were wrong for algo - update HIGH/LOW/LTP/PREV rates on the DLF sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 (daily summary row): F, G, H, I, J
$ws.Range("F7").Value = 678.1
$ws.Range("G7").Value = 694.9
$ws.Range("H7").Value = 670.2
$ws.Range("I7").Value = 692.9
$ws.Range("J7").Value = 675.35

# Row 9: G, H, I
$ws.Range("G9").Value = 685.8
$ws.Range("H9").Value = 670.2
$ws.Range("I9").Value = 671.75

# Row 10: G, H, I
$ws.Range("G10").Value = 681.2
$ws.Range("H10").Value = 671
$ws.Range("I10").Value = 679.9

# Row 11: G, H, I
$ws.Range("G11").Value = 683.2
$ws.Range("H11").Value = 677.85
$ws.Range("I11").Value = 678.45

# Row 12: G, H, I
$ws.Range("G12").Value = 684.4
$ws.Range("H12").Value = 675.4
$ws.Range("I12").Value = 683.75

# Row 13: G, H, I
$ws.Range("G13").Value = 685
$ws.Range("H13").Value = 682.15
$ws.Range("I13").Value = 683.4

# Row 14: G, H, I
$ws.Range("G14").Value = 686.8
$ws.Range("H14").Value = 682.75
$ws.Range("I14").Value = 686.1

# Row 15: G, H, I
$ws.Range("G15").Value = 686.2
$ws.Range("H15").Value = 681.05
$ws.Range("I15").Value = 683.05

# Row 16: G, H, I
$ws.Range("G16").Value = 686.6
$ws.Range("H16").Value = 683.1
$ws.Range("I16").Value = 686.4

# Row 17: G, H, I
$ws.Range("G17").Value = 686.4
$ws.Range("H17").Value = 680.75
$ws.Range("I17").Value = 683.25

# Row 18: G, H, I
$ws.Range("G18").Value = 684.75
$ws.Range("H18").Value = 681
$ws.Range("I18").Value = 682.75

# Row 19: G, H, I
$ws.Range("G19").Value = 684.4
$ws.Range("H19").Value = 681.05
$ws.Range("I19").Value = 683.75

# Row 20: G, H, I
$ws.Range("G20").Value = 690.5
$ws.Range("H20").Value = 683.35
$ws.Range("I20").Value = 690.45

# Row 21: G, H, I
$ws.Range("G21").Value = 694.9
$ws.Range("H21").Value = 689.55
$ws.Range("I21").Value = 694
